$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.402.08'
$ws.Range('E2').Value = '  -0.41%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.846.97'
$ws.Range('E3').Value = '  -0.18%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.99'
$ws.Range('E5').Value = '  -0.99%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6276'
$ws.Range('E6').Value = '  -3.55%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.07%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07606'
$ws.Range('E8').Value = '  +1.37%  '

# Row 9
$ws.Range('E9').Value = '  -0.43%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.45'
$ws.Range('E10').Value = '  -0.11%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '2.099.25'
$ws.Range('E11').Value = '  +13.38%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07716'
$ws.Range('E12').Value = '  +1.07%  '

# Row 13
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6875'
$ws.Range('E13').Value = '  +0.28%  '

# Row 14
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.987'
$ws.Range('E14').Value = '  -0.85%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.91'
$ws.Range('E15').Value = '  -0.88%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.287.63'
$ws.Range('E16').Value = '  +8.54%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000009906'
$ws.Range('E17').Value = '  +4.56%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.162'
$ws.Range('E18').Value = '  +0.68%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '29.659.01'
$ws.Range('E19').Value = '  +0.35%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '231.57'

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.52'
$ws.Range('E21').Value = '  -0.78%  '

# Row 22
$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('E22').Value = '  -0.03%  '

# Row 23
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.607'
$ws.Range('E23').Value = '  -1.17%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9996'
$ws.Range('E24').Value = '  -0.09%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.66'
$ws.Range('E25').Value = '  -1.77%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1389'
$ws.Range('E26').Value = '  -2.13%  '

# Row 27
$ws.Range('E27').Value = '  -0.60%  '

# Row 28
$ws.Range('E28').Value = '  -0.98%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.471'
$ws.Range('E29').Value = '  -1.27%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05806'
$ws.Range('E30').Value = '  -4.43%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.254'
$ws.Range('E31').Value = '  -0.51%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.119'
$ws.Range('E32').Value = '  -0.55%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.018'
$ws.Range('E33').Value = '  -1.32%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.861'
$ws.Range('E34').Value = '  +0.03%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.158'
$ws.Range('E35').Value = '  -2.52%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7187'
$ws.Range('E36').Value = '  -1.00%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.599'
$ws.Range('E37').Value = '  +0.16%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.249.61'
$ws.Range('E38').Value = '  +4.00%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.791'
$ws.Range('E39').Value = '  -0.21%  '

# Row 40
$ws.Range('E40').Value = '  +1.11%  '

# Row 41
$ws.Range('B41').Value = 'RocketPoolETH'
$ws.Range('C41').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.205.50'
$ws.Range('E41').Value = '  +9.40%  '

# Row 42
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9082'
$ws.Range('E42').Value = '  -0.03%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.080'
$ws.Range('E43').Value = '  -2.38%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9990'
$ws.Range('E44').Value = '  -0.17%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.85'
$ws.Range('E45').Value = '  +0.15%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '67.46'
$ws.Range('E46').Value = '  +1.36%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.324'
$ws.Range('E47').Value = '  -1.70%  '

# Row 48
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.162'
$ws.Range('E48').Value = '  +0.71%  '

# Row 49
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000117'
$ws.Range('E49').Value = '  -5.20%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4019'
$ws.Range('E50').Value = '  -0.82%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.697'
$ws.Range('E51').Value = '  +2.59%  '
